$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8, 1).Value = "2025-08-11 17:20"
$ws.Cells.Item(8, 2).Value = "yaren"
$ws.Cells.Item(8, 3).Value = "cvf"
$ws.Cells.Item(8, 4).Value = 35
$ws.Cells.Item(8, 5).Value = 29
$ws.Cells.Item(8, 6).Value = 24
$ws.Cells.Item(8, 7).Value = 23
$ws.Cells.Item(8, 8).Value = 32
$ws.Cells.Item(8, 9).Value = 23.5
$ws.Cells.Item(8, 10).Value = 0.33
$ws.Cells.Item(8, 11).Value = 0.34
$ws.Cells.Item(8, 12).Value = 0.16
$ws.Cells.Item(8, 13).Value = 0.17
$ws.Cells.Item(8, 14).Value = "%32.64"
$ws.Cells.Item(8, 15).Value = "%34.03"
$ws.Cells.Item(8, 16).Value = "%16.32"
$ws.Cells.Item(8, 17).Value = "%17.01"

# Row 9
$ws.Cells.Item(9, 1).Value = "2025-08-12 07:50"
$ws.Cells.Item(9, 2).Value = "yaren"
$ws.Cells.Item(9, 3).Value = "yaren"
$ws.Cells.Item(9, 4).Value = 31
$ws.Cells.Item(9, 5).Value = 33
$ws.Cells.Item(9, 6).Value = 30
$ws.Cells.Item(9, 7).Value = 31
$ws.Cells.Item(9, 8).Value = 32
$ws.Cells.Item(9, 9).Value = 30.5
$ws.Cells.Item(9, 10).Value = 0.42
$ws.Cells.Item(9, 11).Value = 0.24
$ws.Cells.Item(9, 12).Value = 0.21
$ws.Cells.Item(9, 13).Value = 0.12
$ws.Cells.Item(9, 14).Value = "%42.36"
$ws.Cells.Item(9, 15).Value = "%24.31"
$ws.Cells.Item(9, 16).Value = "%21.18"
$ws.Cells.Item(9, 17).Value = "%12.15"

# Row 10
$ws.Cells.Item(10, 1).Value = "2025-08-12 07:54"
$ws.Cells.Item(10, 2).Value = "yaren"
$ws.Cells.Item(10, 3).Value = "yaren"
$ws.Cells.Item(10, 4).Value = 33
$ws.Cells.Item(10, 5).Value = 31
$ws.Cells.Item(10, 6).Value = 27
$ws.Cells.Item(10, 7).Value = 29
$ws.Cells.Item(10, 8).Value = 32
$ws.Cells.Item(10, 9).Value = 28
$ws.Cells.Item(10, 10).Value = 0.39
$ws.Cells.Item(10, 11).Value = 0.28
$ws.Cells.Item(10, 12).Value = 0.19
$ws.Cells.Item(10, 13).Value = 0.14
$ws.Cells.Item(10, 14).Value = "%38.89"
$ws.Cells.Item(10, 15).Value = "%27.78"
$ws.Cells.Item(10, 16).Value = "%19.44"
$ws.Cells.Item(10, 17).Value = "%13.89"

# Row 11
$ws.Cells.Item(11, 1).Value = "2025-08-12 07:55"
$ws.Cells.Item(11, 2).Value = "yaren"
$ws.Cells.Item(11, 3).Value = "yaren"
$ws.Cells.Item(11, 4).Value = 45
$ws.Cells.Item(11, 5).Value = 48
$ws.Cells.Item(11, 6).Value = 45
$ws.Cells.Item(11, 7).Value = 44
$ws.Cells.Item(11, 8).Value = 46.5
$ws.Cells.Item(11, 9).Value = 44.5
$ws.Cells.Item(11, 10).Value = 0.9
$ws.Cells.Item(11, 11).Value = 0.07
$ws.Cells.Item(11, 12).Value = 0.03
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = "%89.81"
$ws.Cells.Item(11, 15).Value = "%7.06"
$ws.Cells.Item(11, 16).Value = "%2.9"
$ws.Cells.Item(11, 17).Value = "%0.23"

# Row 12
$ws.Cells.Item(12, 1).Value = "2025-08-12 07:56"
$ws.Cells.Item(12, 2).Value = "yaren"
$ws.Cells.Item(12, 3).Value = "yaren"
$ws.Cells.Item(12, 4).Value = 33
$ws.Cells.Item(12, 5).Value = 21
$ws.Cells.Item(12, 6).Value = 25
$ws.Cells.Item(12, 7).Value = 20
$ws.Cells.Item(12, 8).Value = 27
$ws.Cells.Item(12, 9).Value = 22.5
$ws.Cells.Item(12, 10).Value = 0.26
$ws.Cells.Item(12, 11).Value = 0.3
$ws.Cells.Item(12, 12).Value = 0.21
$ws.Cells.Item(12, 13).Value = 0.23
$ws.Cells.Item(12, 14).Value = "%26.37"
$ws.Cells.Item(12, 15).Value = "%29.88"
$ws.Cells.Item(12, 16).Value = "%20.51"
$ws.Cells.Item(12, 17).Value = "%23.24"

# Row 13
$ws.Cells.Item(13, 1).Value = "2025-08-12 08:09"
$ws.Cells.Item(13, 2).Value = "yaren"
$ws.Cells.Item(13, 3).Value = "yaren"
$ws.Cells.Item(13, 4).Value = 19
$ws.Cells.Item(13, 5).Value = 25
$ws.Cells.Item(13, 6).Value = 28
$ws.Cells.Item(13, 7).Value = 24
$ws.Cells.Item(13, 8).Value = 22
$ws.Cells.Item(13, 9).Value = 26
$ws.Cells.Item(13, 10).Value = 0.25
$ws.Cells.Item(13, 11).Value = 0.21
$ws.Cells.Item(13, 12).Value = 0.29
$ws.Cells.Item(13, 13).Value = 0.25
$ws.Cells.Item(13, 14).Value = "%24.83"
$ws.Cells.Item(13, 15).Value = "%21.01"
$ws.Cells.Item(13, 16).Value = "%29.34"
$ws.Cells.Item(13, 17).Value = "%24.83"

# Row 14
$ws.Cells.Item(14, 1).Value = "2025-08-12 08:19"
$ws.Cells.Item(14, 2).Value = "asdf"
$ws.Cells.Item(14, 3).Value = "asdf"
$ws.Cells.Item(14, 4).Value = 36
$ws.Cells.Item(14, 5).Value = 23
$ws.Cells.Item(14, 6).Value = 25
$ws.Cells.Item(14, 7).Value = 29
$ws.Cells.Item(14, 8).Value = 29.5
$ws.Cells.Item(14, 9).Value = 27
$ws.Cells.Item(14, 10).Value = 0.35
$ws.Cells.Item(14, 11).Value = 0.27
$ws.Cells.Item(14, 12).Value = 0.22
$ws.Cells.Item(14, 13).Value = 0.17
$ws.Cells.Item(14, 14).Value = "%34.57"
$ws.Cells.Item(14, 15).Value = "%26.89"
$ws.Cells.Item(14, 16).Value = "%21.68"
$ws.Cells.Item(14, 17).Value = "%16.86"
